$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H38").Value = 251.6923
$ws.Range("I38").Value = 251.6923
$ws.Range("K38").Value = 755.0769
$ws.Range("M38").Value = -383.0769
$ws.Range("H43").Value = 3518869.2
$ws.Range("I43").Value = 13161469
$ws.Range("J43").Value = 12469.363
$ws.Range("K43").Value = 13161469
$ws.Range("L43").Value = 12469.363
$ws.Range("M43").Value = -13161400
$ws.Range("N43").Value = -12607.363
$ws.Range("H44").Value = 76999.5
$ws.Range("J44").Value = 76999.5
$ws.Range("L44").Value = 76999.5
$ws.Range("N44").Value = -77923.5
$ws.Range("H64").Value = 5822.769
$ws.Range("I64").Value = 6049.6
$ws.Range("J64").Value = 5066.6665
$ws.Range("K64").Value = 6049.6
$ws.Range("L64").Value = 5066.6665
$ws.Range("M64").Value = -5801.6
$ws.Range("N64").Value = -5562.6665
$ws.Range("H67").Value = 5822.769
$ws.Range("I67").Value = 6049.6
$ws.Range("J67").Value = 5066.6665
$ws.Range("K67").Value = 6049.6
$ws.Range("L67").Value = 5066.6665
$ws.Range("M67").Value = -5191.6
$ws.Range("N67").Value = -6782.6665
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 7132
$ws.Range("I113").Value = 3024.75
$ws.Range("J113").Value = 9870.166999999999
$ws.Range("K113").Value = 3024.75
$ws.Range("L113").Value = 9870.166999999999
$ws.Range("M113").Value = 229.25
$ws.Range("N113").Value = -16378.167
$ws.Range("H132").Value = 5818.1177
$ws.Range("I132").Value = 4647.0347
$ws.Range("K132").Value = 13941.1041
$ws.Range("M132").Value = -11411.1041
$ws.Range("H133").Value = 95593.336
$ws.Range("J133").Value = 95593.336
$ws.Range("L133").Value = 95593.336
$ws.Range("N133").Value = -105713.336
$ws.Range("H137").Value = 2543.8708
$ws.Range("I137").Value = 1847.579
$ws.Range("K137").Value = 5542.737
$ws.Range("M137").Value = -2992.737
$ws.Range("H141").Value = 2207.8635
$ws.Range("I141").Value = 1222.3529
$ws.Range("K141").Value = 3667.0587
$ws.Range("M141").Value = 1512.9413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 12900
$ws.Range("I5").Value = 17066.666
$ws.Range("K5").Value = 17066.666
$ws.Range("M5").Value = -16954.666
$ws.Range("H61").Value = 3643.9211
$ws.Range("I61").Value = 3240.647
$ws.Range("J61").Value = 7071.75
$ws.Range("K61").Value = 3240.647
$ws.Range("L61").Value = 7071.75
$ws.Range("M61").Value = -3028.647
$ws.Range("N61").Value = -7495.75
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 2834.3062
$ws.Range("I132").Value = 1426
$ws.Range("K132").Value = 4278
$ws.Range("M132").Value = -1748
$ws.Range("H136").Value = 3643.9211
$ws.Range("I136").Value = 3240.647
$ws.Range("J136").Value = 7071.75
$ws.Range("K136").Value = 9721.940999999999
$ws.Range("L136").Value = 21215.25
$ws.Range("M136").Value = -7171.940999999999
$ws.Range("N136").Value = -26315.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 12900
$ws.Range("I4").Value = 17066.666
$ws.Range("K4").Value = 17066.666
$ws.Range("M4").Value = -16951.666
$ws.Range("H86").Value = 14070.5
$ws.Range("I86").Value = 26499.25
$ws.Range("J86").Value = 1641.75
$ws.Range("K86").Value = 26499.25
$ws.Range("L86").Value = 1641.75
$ws.Range("M86").Value = -25376.25
$ws.Range("N86").Value = -3887.75
$ws.Range("H89").Value = 14070.5
$ws.Range("I89").Value = 26499.25
$ws.Range("J89").Value = 1641.75
$ws.Range("K89").Value = 132496.25
$ws.Range("L89").Value = 8208.75
$ws.Range("M89").Value = -126880.25
$ws.Range("N89").Value = -19440.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2116.5
$ws.Range("I12").Value = 1099.5
$ws.Range("J12").Value = 2625
$ws.Range("K12").Value = 1099.5
$ws.Range("L12").Value = 2625
$ws.Range("M12").Value = -929.5
$ws.Range("N12").Value = -2965
$ws.Range("H31").Value = 2740.3513
$ws.Range("I31").Value = 2179.6296
$ws.Range("K31").Value = 2179.6296
$ws.Range("M31").Value = -1884.6296
$ws.Range("H34").Value = 2740.3513
$ws.Range("I34").Value = 2179.6296
$ws.Range("K34").Value = 2179.6296
$ws.Range("M34").Value = -1977.6296
$ws.Range("H58").Value = 1992.8948
$ws.Range("J58").Value = 1579.4
$ws.Range("L58").Value = 1579.4
$ws.Range("N58").Value = -1985.4
$ws.Range("H59").Value = 104
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H74").Value = 49999
$ws.Range("J74").Value = 49999
$ws.Range("L74").Value = 49999
$ws.Range("N74").Value = -51747
$ws.Range("H77").Value = 49999
$ws.Range("J77").Value = 49999
$ws.Range("L77").Value = 149997
$ws.Range("N77").Value = -158733
$ws.Range("H118").Value = 88000
$ws.Range("J118").Value = 88000
$ws.Range("L118").Value = 88000
$ws.Range("N118").Value = -91314
$ws.Range("H122").Value = 2117.4348
$ws.Range("J122").Value = 2249.625
$ws.Range("L122").Value = 6748.875
$ws.Range("N122").Value = -11648.875
$ws.Range("H134").Value = 6307.6313
$ws.Range("I134").Value = 6068.25
$ws.Range("K134").Value = 18204.75
$ws.Range("M134").Value = -15669.75
$ws.Range("H136").Value = 1992.8948
$ws.Range("J136").Value = 1579.4
$ws.Range("L136").Value = 4738.200000000001
$ws.Range("N136").Value = -9838.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 197.57143
$ws.Range("I2").Value = 194.33333
$ws.Range("K2").Value = 1165.99998
$ws.Range("M2").Value = -1052.99998
$ws.Range("H7").Value = 249
$ws.Range("I7").Value = 234.14285
$ws.Range("J7").Value = 275
$ws.Range("K7").Value = 702.4285500000001
$ws.Range("L7").Value = 825
$ws.Range("M7").Value = -590.4285500000001
$ws.Range("N7").Value = -1049
$ws.Range("H34").Value = 827.1429000000001
$ws.Range("I34").Value = 178
$ws.Range("J34").Value = 2450
$ws.Range("K34").Value = 534
$ws.Range("L34").Value = 7350
$ws.Range("M34").Value = -450
$ws.Range("N34").Value = -7518
$ws.Range("H39").Value = 1881.5
$ws.Range("J39").Value = 2985
$ws.Range("L39").Value = 8955
$ws.Range("N39").Value = -9543
$ws.Range("H55").Value = 6647.3335
$ws.Range("J55").Value = 18999.5
$ws.Range("L55").Value = 56998.5
$ws.Range("N55").Value = -57352.5
$ws.Range("H107").Value = 3493.389
$ws.Range("I107").Value = 106.5
$ws.Range("J107").Value = 3916.75
$ws.Range("K107").Value = 319.5
$ws.Range("L107").Value = 11750.25
$ws.Range("M107").Value = 1600.5
$ws.Range("N107").Value = -15590.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3167.606
$ws.Range("I122").Value = 2968.5386
$ws.Range("J122").Value = 3907
$ws.Range("K122").Value = 8905.6158
$ws.Range("L122").Value = 11721
$ws.Range("M122").Value = -6455.6158
$ws.Range("N122").Value = -16621

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2767.64
$ws.Range("I22").Value = 893.1
$ws.Range("J22").Value = 4017.3333
$ws.Range("K22").Value = 893.1
$ws.Range("L22").Value = 4017.3333
$ws.Range("M22").Value = -598.1
$ws.Range("N22").Value = -4607.3333
$ws.Range("H27").Value = 2767.64
$ws.Range("I27").Value = 893.1
$ws.Range("J27").Value = 4017.3333
$ws.Range("K27").Value = 893.1
$ws.Range("L27").Value = 4017.3333
$ws.Range("M27").Value = -786.1
$ws.Range("N27").Value = -4231.3333
$ws.Range("H55").Value = 410.16666
$ws.Range("J55").Value = 422
$ws.Range("L55").Value = 422
$ws.Range("N55").Value = -768
$ws.Range("H109").Value = 79000
$ws.Range("J109").Value = 79000
$ws.Range("L109").Value = 79000
$ws.Range("N109").Value = -81774
$ws.Range("H122").Value = 7125.84
$ws.Range("I122").Value = 6377.727
$ws.Range("K122").Value = 19133.181
$ws.Range("M122").Value = -16683.181
$ws.Range("H139").Value = 78562
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 80602.14
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 80602.14
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -90882.14

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 570
$ws.Range("I17").Value = 570
$ws.Range("K17").Value = 570
$ws.Range("M17").Value = -398
$ws.Range("H132").Value = 8874.857
$ws.Range("I132").Value = 8062.75
$ws.Range("K132").Value = 24188.25
$ws.Range("M132").Value = -21658.25
$ws.Range("H139").Value = 77544.28999999999
$ws.Range("I139").Value = 80000
$ws.Range("K139").Value = 80000
$ws.Range("M139").Value = -74860
